$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule Metrics Tracking")
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("G6").Formula = "=H5+1"
$ws.Range("H6").Formula = "=G6+13"
$ws.Range("I6").Value = 14
$ws.Range("K6").Value = "Estimates are generally accurate and on track. "
$ws.Columns("H").ColumnWidth = 24.7
$ws.Activate() | Out-Null
$ws.Range("A7").Select() | Out-Null
